# mentors.xlsx: rename sheet tab + reset the saved view (scroll/zoom/selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Mentors (Made by Luisa)" -> "Mentors"
$ws.Name = "Mentors"

# Make sure this sheet/window is the active one before touching view state
$ws.Activate()

# Scroll the view back to the top-left corner (was parked at A78)
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# Zoom in to 130%
$excel.ActiveWindow.Zoom = 130

# Move the selection from D96 to D2
$ws.Range("D2").Select()
